# Append 6 new applicant rows (187-192) to the admissions list, matching
# the new entries submitted on 2025-07-11 / 2025-07-12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns F (JSHIR), I (Telegram), J (Phone) and K (Sana/date) hold
# digit-only or date-like strings that Excel would otherwise silently
# coerce into numbers / dates. Force those columns to Text up front so
# the values are written verbatim (preserving the leading "+" in J and
# not turning the dates into date serials), mirroring how every other
# row in this sheet stores these columns as plain text.
$ws.Range("F187:F192").NumberFormat = "@"
$ws.Range("I187:J192").NumberFormat = "@"
$ws.Range("K187:K192").NumberFormat = "@"

$rows = @(
    @{ Row = 187; A = "Baxtiyorov Abdulboriy  Baxrom o'g'li"; B = "Yurisprudensiya"; C = "Rus tili"; D = "Kunduzgi"; E = "AD5552063"; F = "50612076540028"; G = "Toshkent shahri"; H = "Shayxontohur tumani"; I = "998974047888"; J = "+998950500778"; K = "2025-07-11" },
    @{ Row = 188; A = "Husanova Nozima Shokir Qizi"; B = "Inson resurslarini boshqarish"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AE2838288"; F = "40409923980013"; G = "Samarqand viloyati"; H = "Nurobod tumani"; I = "998938792997"; J = "+998949381690"; K = "2025-07-11" },
    @{ Row = 189; A = "Sheralieva Madina Nuralievna"; B = "Yurisprudensiya"; C = "Rus tili"; D = "Kunduzgi"; E = "AE1846105"; F = "61612078660029"; G = "Toshkent shahri"; H = "Bektemir tumani"; I = "998881887967"; J = "+998970647478"; K = "2025-07-11" },
    @{ Row = 190; A = "Nurbek Madraimov"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD6081351"; F = "52611076520034"; G = "Toshkent shahri"; H = "Yashnaobod tumani"; I = "998974031380"; J = "+998974031380"; K = "2025-07-11" },
    @{ Row = 191; A = "Tulyaganov Jafarbek Kozimbek Ogli"; B = "Ijtimoiy ish"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD7166340"; F = "51404086580054"; G = "Toshkent shahri"; H = "Yashnaobod tumani"; I = "998901110800"; J = "+998901110800"; K = "2025-07-11" },
    @{ Row = 192; A = "toxtabayeva ruxsora zafar qizi"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD4895079"; F = "60209076790015"; G = "Toshkent viloyati"; H = "Piskent tumani"; I = "998955555207"; J = "+998955555207"; K = "2025-07-12" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value = $r.A
    $ws.Cells.Item($n, 2).Value = $r.B
    $ws.Cells.Item($n, 3).Value = $r.C
    $ws.Cells.Item($n, 4).Value = $r.D
    $ws.Cells.Item($n, 5).Value = $r.E
    $ws.Cells.Item($n, 6).Value = $r.F
    $ws.Cells.Item($n, 7).Value = $r.G
    $ws.Cells.Item($n, 8).Value = $r.H
    $ws.Cells.Item($n, 9).Value = $r.I
    $ws.Cells.Item($n, 10).Value = $r.J
    $ws.Cells.Item($n, 11).Value = $r.K
}
